$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Blue font color used for both new labels (hex 0000FF -> VBA RGB(0,0,255))
$blue = 16711680

# --- Shape "A" (CasellaDiTesto 2) ---
$tbA = $s.Shapes.AddTextbox(1, 31.04755905511811, 73.52377952755906, 11.809527559055118, 29.081259842519685)
$tbA.Name = "CasellaDiTesto 2"
$tbA.Fill.Visible = $false
$tbA.TextFrame.WordWrap = $true
$tbA.TextFrame.AutoSize = 1

$trA = $tbA.TextFrame.TextRange
$trA.Text = "A"
$trA.Font.Color.RGB = $blue

# --- Shape "B" (CasellaDiTesto 3) ---
$tbB = $s.Shapes.AddTextbox(1, 22.476141732283466, 92.38094488188976, 11.809527559055118, 29.081259842519685)
$tbB.Name = "CasellaDiTesto 3"
$tbB.Fill.Visible = $false
$tbB.TextFrame.WordWrap = $true
$tbB.TextFrame.AutoSize = 1

$trB = $tbB.TextFrame.TextRange
$trB.Text = "B"
$trB.Font.Color.RGB = $blue
